$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "51.915.95"
$ws.Range("E2").Value2 = "  +0.29%  "
$ws.Range("D3").Value2 = "2.823.63"
$ws.Range("E3").Value2 = "  +2.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value2 = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "359.97"
$ws.Range("E5").Value2 = "  +8.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.64"
$ws.Range("E6").Value2 = "  -2.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("E7").Value2 = "  +2.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value2 = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value2 = "  +4.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.85"
$ws.Range("E10").Value2 = "  +0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0852"
$ws.Range("E11").Value2 = "  +2.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.91"
$ws.Range("E12").Value2 = "  -0.49%  "
$ws.Range("E13").Value2 = "  +1.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.80"
$ws.Range("E14").Value2 = "  +2.80%  "
$ws.Range("D15").Value2 = "3.277.95"
$ws.Range("E15").Value2 = "  +2.86%  "
$ws.Range("D16").Value2 = "2.841.89"
$ws.Range("E16").Value2 = "  +2.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.895"
$ws.Range("E17").Value2 = "  +1.38%  "
$ws.Range("D18").Value2 = "51.933.11"
$ws.Range("E18").Value2 = "  +0.36%  "
$ws.Range("B19").Value2 = "ImmutableX"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.16"
$ws.Range("E19").Value2 = "  +3.14%  "
$ws.Range("B20").Value2 = "Uniswap"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value2 = "  +6.28%  "
$ws.Range("B21").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.62"
$ws.Range("E21").Value2 = "  +0.92%  "
$ws.Range("D22").Value2 = "0.0₃0979"
$ws.Range("E22").Value2 = "  +1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.76"
$ws.Range("E23").Value2 = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.03"
$ws.Range("E24").Value2 = "  -3.74%  "
$ws.Range("E25").Value2 = "  +6.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.83"
$ws.Range("E26").Value2 = "  +0.23%  "
$ws.Range("E27").Value2 = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.23"
$ws.Range("E28").Value2 = "  -0.14%  "
$ws.Range("E29").Value2 = "  +1.09%  "
$ws.Range("E30").Value2 = "  -0.04%  "
$ws.Range("B31").Value2 = "OKB"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "50.86"
$ws.Range("E31").Value2 = "  +0.75%  "
$ws.Range("B32").Value2 = "InjectiveProtocol"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.03"
$ws.Range("E32").Value2 = "  -2.81%  "
$ws.Range("B33").Value2 = "VeChain"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0442"
$ws.Range("E33").Value2 = "  +27.77%  "
$ws.Range("B34").Value2 = "Filecoin"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.77"
$ws.Range("E34").Value2 = "  +3.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0827"
$ws.Range("E35").Value2 = "  +0.82%  "
$ws.Range("E36").Value2 = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("E37").Value2 = "  +0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.95"
$ws.Range("E38").Value2 = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value2 = "  +2.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.39"
$ws.Range("E40").Value2 = "  -3.62%  "
$ws.Range("B41").Value2 = "EnergySwap"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.61"
$ws.Range("E41").Value2 = "  +1.36%  "
$ws.Range("B42").Value2 = "Stacks"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.56"
$ws.Range("E42").Value2 = "  +4.56%  "
$ws.Range("B43").Value2 = "Monero"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.69"
$ws.Range("E43").Value2 = "  -1.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.115"
$ws.Range("E44").Value2 = "  +1.76%  "
$ws.Range("E45").Value2 = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.35"
$ws.Range("E46").Value2 = "  +0.42%  "
$ws.Range("D47").Value2 = "2.058.21"
$ws.Range("E47").Value2 = "  -2.83%  "
$ws.Range("E48").Value2 = "  +3.51%  "
$ws.Range("B49").Value2 = "RocketPoolETH"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value2 = "3.113.17"
$ws.Range("E49").Value2 = "  +2.94%  "
$ws.Range("B50").Value2 = "SEI"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.935"
$ws.Range("E50").Value2 = "  +6.44%  "
$ws.Range("B51").Value2 = "THORChain"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.60"
$ws.Range("E51").Value2 = "  +0.41%  "
